# "Update gh-pages to output generated at 456a3b4"
#
# The scraper re-ran and picked up one more 漫展 (con) listing, plus a
# small correction to the "想去人数" (interested-count) figure on the
# existing LZ栗子动漫游戏嘉年华 row. Both the "展览" sheet and the
# "全部类型" sheet carry the same table, so the update is applied to both.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # 想去人数 for row 2 (丽水·LZ栗子动漫游戏嘉年华) ticked down from 472 to 471
    $ws.Range("F2").Value = 471

    # Clone row 3's formatting onto the new row 4 first (so the "序号"
    # cell A4 picks up the same bordered/centered style used by A2/A3),
    # then fill in the new event's data.
    $ws.Range("A3").Copy()
    $ws.Range("A4").PasteSpecial(-4122)   # xlPasteFormats

    $ws.Cells.Item(4, 1).Value = 3

    # Force the start-date into B4 as literal text (matching how the
    # existing date-like strings in B2/B3 are stored) rather than letting
    # Excel auto-convert "2024-10-02" into a date serial.
    $dateCell = $ws.Cells.Item(4, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2024-10-02"
    $dateCell.Style = "Normal"

    $ws.Cells.Item(4, 3).Value = "青田·未闻展名国漫嘉年华"
    $ws.Cells.Item(4, 4).Value = "瓯南街道百悦城4幢 西娜君澜大饭店"
    $ws.Cells.Item(4, 5).Value = "2024.10.02 09:00-10.02 17:00"
    $ws.Cells.Item(4, 6).Value = 0
    $ws.Cells.Item(4, 7).Value = 39.9
    $ws.Cells.Item(4, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91328"
    $ws.Cells.Item(4, 9).Value = "//i0.hdslb.com/bfs/openplatform/202408/w8uKtdlg1724147282076.jpeg"
}
